# Append the new commit-log row (row 18) reported by PR #47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = "edit2"
$ws.Cells.Item($row, 2).Value = "riya-morankar"
$ws.Cells.Item($row, 3).Value = "Merged"
$ws.Cells.Item($row, 4).Value = "N/A"

# Leading apostrophe forces the date-looking string to stay literal text
# (matching the sheet's existing "Date" column values) instead of being
# auto-converted into a serial date number; resetting the style afterward
# drops the quote-prefix formatting Excel would otherwise tag the cell with.
$ws.Cells.Item($row, 5).Value = "'2025-06-19"
$ws.Cells.Item($row, 5).Style = "Normal"

$ws.Cells.Item($row, 6).Value = "7d293280560b9ddde7e5a9a9704e8051967a7a1e"
